# Auto-generated COM-interop script to apply the commit diff
$wb = $excel.ActiveWorkbook

$wsInv = $wb.Worksheets.Item("investigaciones")
$wsInv.Range("C11").Value = "Caracterización de la oferta y demanda cultural en horarios no convencionales, con énfasis en públicos, franjas horarias e impactos. Incluye un módulo sobre estímulos e hitos 24/7 del sector.`n"
$wsInv.Range("C15").Value = "La investigación se orienta a la caracterización de públicos y al análisis de percepciones asociadas a los eventos culturales de gran formato de Idartes, con el fin de evaluar los impactos culturales, sociales y económicos de los Festivales al Parque. El estudio analiza la composición y perfiles de asistencia, los niveles de satisfacción, los patrones de consumo, la percepción de marca y los aportes de los festivales al posicionamiento de Bogotá como ciudad de grandes eventos. Los resultados constituyen insumos estratégicos para la toma de decisiones, el fortalecimiento de la oferta cultural y la proyección de la ciudad a nivel nacional e internacional."
$wsInv.Range("B21").Value = "Medición de eventos culturales de gran formato del IDPC"
$wsInv.Range("C21").Value = "La investigación se orienta a la caracterización de públicos y al análisis de percepciones asociadas a los eventos culturales de gran formato del IDPC, con el fin de evaluar los impactos culturales, sociales y económicos de la Noche de Museos, la Noche Iberoamericana de Museos y el Festival de Patrimonios en Ruana. El estudio analiza la composición y perfiles de asistencia, los niveles de satisfacción, los patrones de consumo, la percepción de marca y los aportes de los festivales al posicionamiento de Bogotá como ciudad de grandes eventos. Los resultados constituyen insumos estratégicos para la toma de decisiones, el fortalecimiento de la oferta cultural y la proyección de la ciudad a nivel nacional e internacional."
$wsInv.Range("C26").Value = "La investigación se orienta a la caracterización de públicos y al análisis de percepciones asociadas a los eventos culturales de gran formato organizados por la SCRD con el fin de evaluar los impactos culturales, sociales y económicos del Concurso Internacional de Violín, la Feria de Navidad, el Encuentro de Cultura en Iberoamérica y la Bienal de Arte y Ciudad. El estudio analiza la composición y perfiles de asistencia, los niveles de satisfacción, los patrones de consumo, la percepción de marca y los aportes de los festivales al posicionamiento de Bogotá como ciudad de grandes eventos. Los resultados constituyen insumos estratégicos para la toma de decisiones, el fortalecimiento de la oferta cultural y la proyección de la ciudad a nivel nacional e internacional."
$wsInv.Range("C27").Value = "La investigación se orienta a la caracterización de públicos y al análisis de percepciones asociadas a los eventos culturales de gran formato del IDRD, con el fin de evaluar los impactos culturales, sociales y económicos del Festival de Verano. El estudio analiza la composición y perfiles de asistencia, los niveles de satisfacción, los patrones de consumo, la percepción de marca y los aportes de los festivales al posicionamiento de Bogotá como ciudad de grandes eventos. Los resultados constituyen insumos estratégicos para la toma de decisiones, el fortalecimiento de la oferta cultural y la proyección de la ciudad a nivel nacional e internacional.`n"

$wsProd = $wb.Worksheets.Item("productos")
$wsProd.Range("D550").Value = "Sí"
$wsProd.Range("E578").Value = "https://drive.google.com/file/d/1zSXYH6cCKPEoIcuur0PuPncZ0IaOLo8E/view?usp=drive_link"
$wsProd.Range("B582").Value = "Instrumento recolección"
$wsProd.Range("C582").Value = "Formularios Festival Monumentum 2025"
$wsProd.Range("D582").Value = "Sí"
$wsProd.Range("E582").Value = "https://drive.google.com/file/d/1f8Q-GMDfGsdYbhOeIsTIlHRfCnQpw8vh/view?usp=sharing"
$wsProd.Range("B594").Value = "Instrumento recolección"
$wsProd.Range("C594").Value = "Formularios Premio Luis Caballero 2025"
$wsProd.Range("D594").Value = "Sí"
$wsProd.Range("E594").Value = "https://drive.google.com/file/d/1tA-OMX6rjr4wNMQ0GIEXPrZ_v1L66jet/view?usp=sharing"
$wsProd.Range("A597").Value = 108
$wsProd.Range("B597").Value = "Presentación"
$wsProd.Range("C597").Value = "Presentación resultados Festivales al Parque 2025"
$wsProd.Range("D597").Value = "Sí"
$wsProd.Range("E597").Value = "https://drive.google.com/file/d/1ju-2pVTUW2-guIgxcUm7FLAnb1yNjvGX/view?usp=sharing"
$wsProd.Range("E605").Value = "https://drive.google.com/drive/folders/1Lcu0iYLFoiXHMX8f8ivPgJPF5XDn-IYx?usp=sharing"
$wsProd.Range("E606").Value = "https://drive.google.com/drive/folders/1rccxA6SFUSsKbjgRKqWnvEbok0SpWS9t?usp=sharing"
$wsProd.Range("E607").Value = "https://drive.google.com/drive/folders/1BchZfxR7zuYYbkgIH8qE8fdgQ1y1EYZk?usp=sharing"
$wsProd.Range("E608").Value = "https://drive.google.com/drive/folders/1HBmOWH-hQDqwj0mDiTDqlVdSUVcQkhYt?usp=sharing"
$wsProd.Range("E609").Value = "https://drive.google.com/drive/folders/1SSz2sXnITpz_lT4g5-n2AR3M8Ed3LsnM?usp=sharing"
$wsProd.Range("E616").Value = "https://drive.google.com/drive/folders/14CmZAQrG8eHOh7xDTf8ypS7TweGA5Df_?usp=sharing"
$wsProd.Range("E617").Value = "https://drive.google.com/drive/folders/19PXuTzH3C_Nfl5Gj8LO1_XHUgBfjb4b4?usp=sharing"
$wsProd.Range("E618").Value = "https://drive.google.com/drive/folders/1DwWfpUiZRNgalk9jBxS-hFLcAWPqdcap?usp=sharing"
$wsProd.Range("B619").Value = "Carpeta archivos"
$wsProd.Range("C619").Value = "Productos finales Festival Patrimonios en Ruana 2025"
$wsProd.Range("E619").Value = "https://drive.google.com/drive/folders/11bbVLTsTRBTvuwOBfs-eO66V-l-ZCh3_?usp=sharing"
$wsProd.Range("B620").Value = "Carpeta archivos"
$wsProd.Range("C620").Value = "Productos finales Noche de Museos 2025"
$wsProd.Range("E620").Value = "https://drive.google.com/drive/folders/10FyRwNotpVpQk4mvg9tOUyDOFYRDlXA3?usp=sharing"
$wsProd.Range("B621").Value = ""
$wsProd.Range("C621").Value = ""
$wsProd.Range("D621").Value = ""
$wsProd.Range("E621").Value = ""
$wsProd.Range("B622").Value = ""
$wsProd.Range("C622").Value = ""
$wsProd.Range("D622").Value = ""
$wsProd.Range("E622").Value = ""
$wsProd.Range("E670").Value = "https://drive.google.com/file/d/1-uaRSsIb95aqNSXjLDLMJG0a-OjealWM/view?usp=sharing"
$wsProd.Range("B671").Value = "Instrumento recolección"
$wsProd.Range("C671").Value = "Formularios Festival de Verano 2025"
$wsProd.Range("D671").Value = "Sí"
$wsProd.Range("E671").Value = "https://drive.google.com/file/d/1mDfyDNhmr4Y0uawfGhle1_DlUiM9H9oU/view?usp=sharing"
$wsProd.Range("B676").Value = "Carpeta archivos"
$wsProd.Range("C676").Value = "Productos finales Concurso Internacional de Violín 2025"
$wsProd.Range("E676").Value = "https://drive.google.com/drive/folders/114QuIa8ZzjC_GNjRUVggtFVTN3y_zOAL?usp=sharing"
$wsProd.Range("B677").Value = "Carpeta archivos"
$wsProd.Range("C677").Value = "Productos finales Bienal Internacional de Arte y Ciudad BOG25"
$wsProd.Range("E677").Value = "https://drive.google.com/drive/folders/17Bi6a6nKxYmJNqTEWtIH6lYiXMxIECYi?usp=sharing"

$wsHal = $wb.Worksheets.Item("hallazgos")
$wsHal.Rows.Item(311).Insert()
$wsHal.Rows.Item(312).Insert()
$wsHal.Rows.Item(313).Insert()
$wsHal.Rows.Item(314).Insert()
$wsHal.Range("A311").Value = 108
$wsHal.Range("B311").Value = 1
$wsHal.Range("C311").Value = "Experiencia y calidad artística como eje de satisfacción"
$wsHal.Range("D311").Value = "En los ocho Festivales al Parque 2025, el balance general converge en una valoración muy favorable de la experiencia, donde el disfrute del espectáculo, la calidad artística y la programación aparecen como el principal motor de satisfacción y como el rasgo que sostiene la reputación de los eventos. En conjunto, los eventos se consolidan como plataformas culturales del Distrito que ofrecen experiencias significativas en espacio público, con una percepción positiva sobre la organización y logística, reforzando su carácter emblemático dentro de la vida cultural de Bogotá D.C."
$wsHal.Range("A312").Value = 108
$wsHal.Range("B312").Value = 2
$wsHal.Range("C312").Value = "Renovación y fidelización de públicos"
$wsHal.Range("D312").Value = "Los resultados muestran una dinámica consistente de públicos mixtos con personas que asisten por primera vez y que conviven con asistentes recurrentes, lo que sugiere simultáneamente capacidad de ampliar audiencias y de fidelizar comunidades culturales ya formadas. Esta combinación refuerza el rol de los festivales como puerta de entrada a la oferta cultural distrital y, al mismo tiempo, como un circuito de permanencia donde se sostienen trayectorias de participación. Además, se reconoce que la experiencia del evento tiende a motivar a los asistentes a explorar otras actividades culturales de la ciudad."
$wsHal.Range("A313").Value = 108
$wsHal.Range("B313").Value = 3
$wsHal.Range("C313").Value = "Comunidad, convivencia e inclusión"
$wsHal.Range("D313").Value = "De manera transversal, los festivales son percibidos como espacios de encuentro ciudadano que favorecen la convivencia y el reconocimiento mutuo entre públicos diversos. Sumado al componente musical, resalta la capacidad de los Festivales al Parque para activar vínculos sociales, fortalecer sentidos de pertenencia y producir experiencias colectivas en torno a prácticas culturales compartidas. "
$wsHal.Range("A314").Value = 108
$wsHal.Range("B314").Value = 4
$wsHal.Range("C314").Value = "Identidad y orgullo por Bogotá D.C."
$wsHal.Range("D314").Value = "Se identifica una relación entre los festivales y el fortalecimiento de la identidad cultural, tanto por el reconocimiento musical y de expresiones artísticas, como por la apropiación simbólica del espacio público, lo cual se traduce en percepciones favorables de orgullo por la ciudad y por su oferta cultural. De manera complementaria, se observan valoraciones menos homogéneas en lo que respecta a la confianza institucional y algunos componentes de la experiencia urbana, lo que sugiere un campo claro para fortalecer el vínculo entre experiencia cultural e institucionalidad."
